$p = $ppt.ActivePresentation
$s = $p.Slides.Item(71)
$sh = $s.Shapes.Item(2)
$tf2 = $sh.TextFrame2
$tr2 = $tf2.TextRange
try {
    $tr2.InsertXML("<m:oMath xmlns:m='http://schemas.openxmlformats.org/officeDocument/2006/math'><m:r><m:t>x</m:t></m:r></m:oMath>")
    Write-Host "OK"
} catch {
    Write-Host "ERR: $_"
}
